$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 117, shifting the existing
# rows 117-203 down to 119-205 (this also grows the used range to R205).
$ws.Rows("117:118").Insert()

# New weekly data point (date 44574) inserted at the top of this
# Vega Monumental Concepcion / Acelga block: one row for "Primera"
# quality and one for "Segunda" quality.
$row117 = @(11, "Vega Monumental Concepción", "Bíobío", 44574, 8, 100112009, "Acelga", "Sin especificar", "Primera", 200, 600, 700, 650, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 650, 1, "Hortaliza")
$row118 = @(11, "Vega Monumental Concepción", "Bíobío", 44574, 8, 100112009, "Acelga", "Sin especificar", "Segunda", 100, 500, 500, 500, "`$/atado 0,5 a 1 kilo", "Región de Ñuble", 500, 1, "Hortaliza")

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(117, $col).Value = $row117[$col - 1]
    $ws.Cells.Item(118, $col).Value = $row118[$col - 1]
}
